# Rename the descr/AlternativeText references on the pictures of slides
# 16-21 from "cell-16-output-*.png" to "cell-10-output-*.png".
$p = $ppt.ActivePresentation

for ($i = 16; $i -le 21; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.AlternativeText -like "*cell-16-output-*") {
            $sh.AlternativeText = $sh.AlternativeText -replace "cell-16-output-", "cell-10-output-"
        }
    }
}
